# Generate Report for Handoff
#
# The previous handoff round used source id "46202592-973e-49bb-ba9d-3f391884eccf"
# and a particular xliff content hash / set of timestamps. This run
# regenerates the report for the new handoff round, which has a new
# source id "7863b9b3-b309-4b98-8f54-0a6c72c729a7", a new xliff content
# hash, and refreshed handoff timestamps.

$wb = $excel.ActiveWorkbook

$oldId  = "46202592-973e-49bb-ba9d-3f391884eccf"
$newId  = "7863b9b3-b309-4b98-8f54-0a6c72c729a7"

$oldHash = "6944f992e2350c25c5ad34be727e5f596b46c276"
$newHash = "264d8e8bb73fcef1ed9c8c1f0a4b6768e8947973"

$newLatestHoDate   = "2016-09-04 05:04:33"
$newHandoffDateZh  = "2016-09-04 05:04:28"

$mdFile      = "$newId.md"
$mdDisplay   = "e2e\$newId.md"
$zhXlf       = "$newId.$newHash.zh-cn.xlf"
$deXlf       = "$newId.$newHash.de-de.xlf"

# Hyperlink target addresses are unchanged by this edit (still pointing at
# the old commit/id in the repo); only the displayed text is refreshed.
$hyperlinkAddress = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/01ec8d788a9200519d69ad0c6e7d9f88298b0954/e2e/$oldId.md"

function Update-HyperlinkDisplay($ws, $cellRef, $displayText) {
    $rng = $ws.Range($cellRef)
    $rng.Hyperlinks.Delete()
    $rng.Hyperlinks.Add($rng, $hyperlinkAddress, [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, $displayText) | Out-Null
}

# ---- Overview sheet ----
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("A2").Value = $mdFile
$wsOverview.Range("B2").Value = $mdDisplay
$wsOverview.Range("G2").Value = $newLatestHoDate
Update-HyperlinkDisplay $wsOverview "B2" $mdDisplay

# ---- zh-cn sheet ----
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("A2").Value = $mdFile
$wsZh.Range("G2").Value = $zhXlf
$wsZh.Range("H2").Value = $newHandoffDateZh
Update-HyperlinkDisplay $wsZh "A2" $mdFile

# ---- de-de sheet ----
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("A2").Value = $mdFile
$wsDe.Range("G2").Value = $deXlf
$wsDe.Range("H2").Value = $newLatestHoDate
Update-HyperlinkDisplay $wsDe "A2" $mdFile
